# Add new quarter column BB (EQUIPMENT eval) and new row 84 to the
# quarterly real-time data matrix on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BB ---------------------------------------------------
# BB1: new vintage date header, same style/format as BA1 (date cell)
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

# BB2:BB81 duplicate the values already present in BA2:BA81 (values only,
# no style - matches BA2:BA81 which carry no explicit style either)
$ws.Range("BA2:BA81").Copy()
$ws.Range("BB2:BB81").PasteSpecial(-4163)

# BB82 and BB83 carry revised (re-estimated) values, different from BA82/BA83
$ws.Range("BB82").Value = -0.5
$ws.Range("BB83").Value = 0.4

# --- New row 84 --------------------------------------------------------
# A84: new observation date, same style as the rest of column A (date format)
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("A84").Value = 45884

# BB84: first (and only, for now) estimate for the new observation
$ws.Range("BB84").Value = 0.2
